$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data for rows 2-21 (A:T), reflecting the refreshed
# NATMI ligand-receptor pair computation ("Natmi following Dr Hou advice").
# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E Ligand-expressing cells, F Ligand detection rate, G Ligand average expression value,
# H Ligand total expression value, I Ligand derived specificity (avg), J Ligand derived specificity (total),
# K Receptor-expressing cells, L Receptor detection rate, M Receptor average expression value,
# N Receptor total expression value, O Receptor derived specificity (avg), P Receptor derived specificity (total),
# Q Edge average expression weight, R Edge total expression weight,
# S Edge average expression derived specificity, T Edge total expression derived specificity
$data = @(
    @("ECs","Gnai2","Tbxa2r","ECs","3","1","151.7260716666667","455.178215","0.2700739458961593","0.2783366498663096","3","1","1.558708666666667","4.676126","0.3590038617594454","0.3716334638810379","236.4967428661211","2128.47068579509","0.0969575895373327","0.1034392133148602"),
    @("ECs","Gnai2","Tbxa2r","FAPs","3","1","151.7260716666667","455.178215","0.2700739458961593","0.2783366498663096","3","1","2.177429","6.532287","0.5015083552327337","0.5191512044104615","330.3727485030784","2973.354736527705","0.1354443403975972","0.1444988070096675"),
    @("ECs","Gnai2","Tbxa2r","M2","3","1","151.7260716666667","455.178215","0.2700739458961593","0.2783366498663096","2","0.6666666666666666","0.162971","0.488913","0.03753569836749389","0.0388561881622672","24.72694962558834","222.542546630295","0.0101374141700771","0.01081510123966041"),
    @("ECs","Gnai2","Tbxa2r","sCs","3","1","151.7260716666667","455.178215","0.2700739458961593","0.2783366498663096","2","1","0.4426515","0.885303","0.1019520846403269","0.07035914354623346","67.1617732123575","402.970639274145","0.02753460179115231","0.0195835283021214"),
    @("FAPs","Gnai2","Tbxa2r","ECs","3","1","82.24887099999999","246.746613","0.146403824289839","0.150882936320401","3","1","1.558708666666667","4.676126","0.3590038617594454","0.3716334638810379","128.2020280512487","1153.818252461238","0.05255953829640348","0.05607314826529267"),
    @("FAPs","Gnai2","Tbxa2r","FAPs","3","1","82.24887099999999","246.746613","0.146403824289839","0.150882936320401","3","1","2.177429","6.532287","0.5015083552327337","0.5191512044104615","179.091076932659","1611.819692393931","0.0734227411193793","0.07833105811572313"),
    @("FAPs","Gnai2","Tbxa2r","M2","3","1","82.24887099999999","246.746613","0.146403824289839","0.150882936320401","2","0.6666666666666666","0.162971","0.488913","0.03753569836749389","0.0388561881622672","13.404180755741","120.637626801669","0.005495369788390971","0.00586273576414088"),
    @("FAPs","Gnai2","Tbxa2r","sCs","3","1","82.24887099999999","246.746613","0.146403824289839","0.150882936320401","2","1","0.4426515","0.885303","0.1019520846403269","0.07035914354623346","36.4075861214565","218.445516728739","0.01492617508566521","0.01061599417524429"),
    @("M1","Gnai2","Tbxa2r","ECs","3","1","123.444321","370.332963","0.2197321429647646","0.2264546783208506","3","1","1.558708666666667","4.676126","0.3590038617594454","0.3716334638810379","192.413732993482","1731.723596941338","0.07888468787702904","0.08415813651644387"),
    @("M1","Gnai2","Tbxa2r","FAPs","3","1","123.444321","370.332963","0.2197321429647646","0.2264546783208506","3","1","2.177429","6.532287","0.5015083552327337","0.5191512044104615","268.791244430709","2419.121199876381","0.110197505610023","0.1175642189946532"),
    @("M1","Gnai2","Tbxa2r","M2","3","1","123.444321","370.332963","0.2197321429647646","0.2264546783208506","2","0.6666666666666666","0.162971","0.488913","0.03753569836749389","0.0388561881622672","20.117844437691","181.060599939219","0.008247799439968449","0.008799165591060662"),
    @("M1","Gnai2","Tbxa2r","sCs","3","1","123.444321","370.332963","0.2197321429647646","0.2264546783208506","2","1","0.4426515","0.885303","0.1019520846403269","0.07035914354623346","54.6428138571315","327.856883142789","0.0224021500377441","0.01593315721869285"),
    @("M2","Gnai2","Tbxa2r","ECs","3","1","154.3429766666667","463.02893","0.2747320633285943","0.2831372788071194","3","1","1.558708666666667","4.676126","0.3590038617594454","0.3716334638810379","240.5757353694644","2165.18161832518","0.09862987168410586","0.105223287676941"),
    @("M2","Gnai2","Tbxa2r","FAPs","3","1","154.3429766666667","463.02893","0.2747320633285943","0.2831372788071194","3","1","2.177429","6.532287","0.5015083552327337","0.5191512044104615","336.0708733403233","3024.63786006291","0.1377804252096186","0.1469910593062167"),
    @("M2","Gnai2","Tbxa2r","M2","3","1","154.3429766666667","463.02893","0.2747320633285943","0.2831372788071194","2","0.6666666666666666","0.162971","0.488913","0.03753569836749389","0.0388561881622672","25.15342925034333","226.38086325309","0.01031225986098134","0.01100163538108174"),
    @("M2","Gnai2","Tbxa2r","sCs","3","1","154.3429766666667","463.02893","0.2747320633285943","0.2831372788071194","2","1","0.4426515","0.885303","0.1019520846403269","0.07035914354623346","68.32015013596499","409.92090081579","0.0280095065738885","0.01992129644288004"),
    @("sCs","Gnai2","Tbxa2r","ECs","2","1","50.0323125","100.064625","0.08905802352064279","0.06118845668531954","3","1","1.558708666666667","4.676126","0.3590038617594454","0.3716334638810379","77.985799107125","467.9147946427501","0.03197217436457428","0.02273967810750015"),
    @("sCs","Gnai2","Tbxa2r","FAPs","2","1","50.0323125","100.064625","0.08905802352064279","0.06118845668531954","3","1","2.177429","6.532287","0.5015083552327337","0.5191512044104615","108.9418081745625","653.650849047375","0.04466334289611568","0.03176606098420099"),
    @("sCs","Gnai2","Tbxa2r","M2","2","1","50.0323125","100.064625","0.08905802352064279","0.06118845668531954","2","0.6666666666666666","0.162971","0.488913","0.03753569836749389","0.0388561881622672","8.153816000437502","48.92289600262501","0.003342855108076024","0.002377550186323513"),
    @("sCs","Gnai2","Tbxa2r","sCs","2","1","50.0323125","100.064625","0.08905802352064279","0.06118845668531954","2","1","0.4426515","0.885303","0.1019520846403269","0.07035914354623346","22.14687817659375","88.58751270637499","0.0090796511518768","0.004305167407294886")
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $rowVals[$j]
    }
}
